$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Start/End timestamps and level values with the finished,
# more precise date/time values from the completed date processing.

$ws.Range("C3").Value = "11-05-202110:40"
$ws.Range("D3").Value = "12-05-202110:50"
$ws.Range("E3").Value = "low"

$ws.Range("C4").Value = "11-05-202111:00"
$ws.Range("D4").Value = "12-05-202111:10"
$ws.Range("E4").Value = "low"

$ws.Range("C5").Value = "31-05-202115:00"
$ws.Range("D5").Value = "01-06-202115:00"
$ws.Range("E5").Value = "High"

$ws.Range("C6").Value = "31-05-202115:30"
$ws.Range("D6").Value = "01-06-202115:30"
$ws.Range("E6").Value = "High"

$ws.Range("C7").Value = "27-06-202113:30"
$ws.Range("D7").Value = "28-06-202112:45"
$ws.Range("E7").Value = "Low"

$ws.Range("C8").Value = "27-06-202113:30"
$ws.Range("D8").Value = "28-06-202113:15"
$ws.Range("E8").Value = "Low"

$ws.Range("C9").Value = "18-07-202118:00"
$ws.Range("D9").Value = "29-07-202118:00"
$ws.Range("E9").Value = "High"

$ws.Range("C10").Value = "18-07-202118:00"
$ws.Range("D10").Value = "29-07-202118:00"
$ws.Range("E10").Value = "High"

$ws.Range("C11").Value = "11-08-202113:00"
$ws.Range("D11").Value = "12-08-202113:40"
$ws.Range("E11").Value = "Low"

$ws.Range("C12").Value = "11-08-202113:40"
$ws.Range("D12").Value = "12-08-202114:00"
$ws.Range("E12").Value = "Low"

$ws.Range("C13").Value = "08-09-202113:20"
$ws.Range("D13").Value = "09-09-202114:00"
$ws.Range("E13").Value = "High"

$ws.Range("C14").Value = "08-09-202113:40"
$ws.Range("D14").Value = "09-09-202114:00"
$ws.Range("E14").Value = "High"

$ws.Range("D13").Select()
